# Update "Poker - Person Info" worksheet with latest data (up to Nov 2023).
# Only the NoOfSessions (D), PointsBonusTotal (I) and AvgSessPoints (J) columns
# change for a handful of players; everything else stays the same.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 (Jon)
$ws.Range("D4").Value = 184
$ws.Range("I4").Value = 666

# Row 5 (Maisy)
$ws.Range("D5").Value = 112
$ws.Range("I5").Value = 420
$ws.Range("J5").Value = 3.75

# Row 6 (Mark)
$ws.Range("D6").Value = 131
$ws.Range("I6").Value = 470
$ws.Range("J6").Value = 3.59

# Row 7 (Matt)
$ws.Range("D7").Value = 179
$ws.Range("I7").Value = 685

# Row 8 (Pepe)
$ws.Range("D8").Value = 94
$ws.Range("I8").Value = 348
$ws.Range("J8").Value = 3.7

# Row 9 (Prashant)
$ws.Range("D9").Value = 29
$ws.Range("I9").Value = 111
$ws.Range("J9").Value = 3.83

# Row 10 (Richard)
$ws.Range("D10").Value = 131
$ws.Range("I10").Value = 546
$ws.Range("J10").Value = 4.17

# Row 12 (Andy)
$ws.Range("D12").Value = 81
$ws.Range("J12").Value = 3.38

# Row 13 (Anthony)
$ws.Range("D13").Value = 188
$ws.Range("I13").Value = 777
$ws.Range("J13").Value = 4.13

# Row 14 (Bob)
$ws.Range("D14").Value = 112
$ws.Range("I14").Value = 441
$ws.Range("J14").Value = 3.94
